$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("equipment")
$ws.Columns.Item(10).Insert()
Write-Host "inserted"
